# Applies the "Next button fixes and text case and powword files separated" edit
# to the "testcase" worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testcase")

# --- Update existing locator/value cells (rows 5-9) -----------------------
# G5: new locator text for the Case table body
$ws.Range("G5").Value = "//div[ contains(text(),'Case')]//parent::span//parent::th//parent::tr//parent::thead//parent::table/tbody"

# G6: new locator text for the Case table header
$ws.Range("G6").Value = "//div[ contains(text(),'Case')]//parent::span//parent::th//parent::tr//parent::thead//parent::table/thead"

# G7: switch Next-button locator back to the hard-coded xpath
$ws.Range("G7").Value = "//*[@id=`"root`"]/div[3]/div/div[4]/div[2]/div/table/tfoot/tr/td[2]/div/div[3]/button[2]"

# Row 8: new "rowcount" / "G_rowcount" step, with G8 stored as text "19"
$ws.Range("C8").Value = "rowcount"
$ws.Range("D8").Value = "G_rowcount"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "19"
$ws.Range("J8").Value = "StoreGlobal"

# Row 9: swap the DB-connect step out for a webdata step
$ws.Range("I9").Value = "na"
$ws.Range("J9").Value = "webdata"

# --- Insert a new row 10 (old row 10 shifts down to row 11) ---------------
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the former "Dbconnect" step data
$ws.Range("A10").Value = "na"
$ws.Range("B10").Value = "na"
$ws.Range("C10").Value = "na"
$ws.Range("D10").Value = "na"
$ws.Range("E10").Value = "na"
$ws.Range("F10").Value = "na"
$ws.Range("G10").Value = "na"
$ws.Range("H10").Value = "na"
$ws.Range("I10").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WITH DISTINCT c AS c, p, s, demo, diag RETURN c.case_id AS `Case ID`, s.clinical_study_designation AS `Study Code`, s.clinical_study_type AS `Study Type`, demo.breed AS Breed, diag.disease_term AS Diagnosis, diag.stage_of_disease AS `Stage of Disease`,demo.patient_age_at_enrollment AS Age, demo.sex AS Sex, demo.neutered_indicator AS `Neutered Status`'
$ws.Range("J10").Value = "Dbconnect"
$ws.Range("K10").Value = "Y"

# --- Update the view's active selection ------------------------------------
$ws.Activate()
$ws.Range("G8").Select()
